$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "24.820.30"
$ws.Range("E2").Value2 = "  +1.25%  "
$ws.Range("D3").Value2 = "1.701.59"
$ws.Range("E3").Value2 = "  +0.61%  "
$ws.Range("D4").Value2 = "'0.9995"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value2 = "  -0.36%  "
$ws.Range("D5").Value2 = "'314.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  +0.07%  "
$ws.Range("E6").Value2 = "  -0.19%  "
$ws.Range("E7").Value2 = "  +2.41%  "
$ws.Range("D8").Value2 = "'0.4038"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value2 = "  +0.00%  "
$ws.Range("B9").Value2 = "Polygon"
$ws.Range("C9").Value2 = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D9").Value2 = "'1.472"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value2 = "  -1.49%  "
$ws.Range("B10").Value2 = "BinanceUSD"
$ws.Range("C10").Value2 = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D10").Value2 = "'1.001"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value2 = "  -0.45%  "
$ws.Range("E11").Value2 = "  +2.00%  "
$ws.Range("D12").Value2 = "'0.08801"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value2 = "  +0.41%  "
$ws.Range("D13").Value2 = "'26.09"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value2 = "  +3.48%  "
$ws.Range("D14").Value2 = "'7.564"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = "  +0.42%  "
$ws.Range("D15").Value2 = "'7.994"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value2 = "  +0.04%  "
$ws.Range("D16").Value2 = "'0.00001349"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value2 = "  -0.31%  "
$ws.Range("D17").Value2 = "1.677.81"
$ws.Range("E17").Value2 = "  -0.50%  "
$ws.Range("D18").Value2 = "'95.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value2 = "  -2.91%  "
$ws.Range("D19").Value2 = "'0.07181"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = "  +1.00%  "
$ws.Range("D20").Value2 = "'20.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "  +4.22%  "
$ws.Range("D21").Value2 = "'7.347"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = "  +0.87%  "
$ws.Range("D22").Value2 = "'1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = "  -0.23%  "
$ws.Range("D23").Value2 = "'14.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = "  +0.65%  "
$ws.Range("D24").Value2 = "24.780.25"
$ws.Range("E24").Value2 = "  +1.10%  "
$ws.Range("D25").Value2 = "'2.368"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value2 = "  +0.73%  "
$ws.Range("D26").Value2 = "'2.916"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = "  -2.09%  "
$ws.Range("D27").Value2 = "'23.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value2 = "  +1.75%  "
$ws.Range("D28").Value2 = "'6.158"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = "  +17.82%  "
$ws.Range("D29").Value2 = "'161.73"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value2 = "  -0.45%  "
$ws.Range("D30").Value2 = "'144.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = "  +5.00%  "
$ws.Range("D31").Value2 = "'8.404"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value2 = "  -4.53%  "
$ws.Range("D32").Value2 = "'2.362"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value2 = "  +20.54%  "
$ws.Range("D33").Value2 = "1.895.02"
$ws.Range("E33").Value2 = "  +1.23%  "
$ws.Range("D34").Value2 = "'0.08648"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value2 = "  -2.29%  "
$ws.Range("D35").Value2 = "'7.326"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value2 = "  -1.76%  "
$ws.Range("D36").Value2 = "'0.03174"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = "  +8.69%  "
$ws.Range("D37").Value2 = "'1.036"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value2 = "  -0.24%  "
$ws.Range("D38").Value2 = "'0.2834"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = "  +0.82%  "
$ws.Range("B39").Value2 = "Stellar"
$ws.Range("C39").Value2 = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").Value2 = "'0.09473"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value2 = "  +3.70%  "
$ws.Range("D40").Value2 = "'10.76"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = "  -0.38%  "
$ws.Range("B41").Value2 = "TheSandbox"
$ws.Range("C41").Value2 = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value2 = "'0.8309"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = "  +4.66%  "
$ws.Range("D42").Value2 = "'14.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value2 = "  -0.31%  "
$ws.Range("E43").Value2 = "  +1.60%  "
$ws.Range("E44").Value2 = "  +6.32%  "
$ws.Range("D45").Value2 = "'2.704"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value2 = "  +3.19%  "
$ws.Range("D46").Value2 = "'0.7436"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = "  +2.47%  "
$ws.Range("D47").Value2 = "'4.215"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = "  +0.24%  "
$ws.Range("D48").Value2 = "'1.386"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "  +2.38%  "
$ws.Range("E49").Value2 = "  -0.07%  "
$ws.Range("D50").Value2 = "'0.08381"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = "  +4.40%  "
$ws.Range("D51").Value2 = "'139.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value2 = "  +1.02%  "
